$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# "About" sheet: the version banner (A2) and the recommended-citation text (A6)
# both embed the build timestamp.
$wsAbout = $wb.Worksheets.Item("About")
foreach ($addr in @("A2", "A6")) {
    $cell = $wsAbout.Range($addr)
    $v = $cell.Value2
    if ($v -ne $null -and $v -is [string] -and $v.Contains($oldStamp)) {
        $cell.Value2 = $v.Replace($oldStamp, $newStamp)
    }
}

# "Boundaries and methane sources" sheet: the build_version column (S2:S12)
# repeats the same version string for every data row.
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 12; $r++) {
    $cell = $wsData.Range("S$r")
    $v = $cell.Value2
    if ($v -ne $null -and $v -is [string] -and $v.Contains($oldStamp)) {
        $cell.Value2 = $v.Replace($oldStamp, $newStamp)
    }
}
